$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) - "Block Trial" column inserted between Trial and Stimulus
$ws.Cells.Item(1, 1).Value = "Trial"
$ws.Cells.Item(1, 2).Value = "Block Trial"
$ws.Cells.Item(1, 3).Value = "Stimulus"
$ws.Cells.Item(1, 4).Value = "Response"
$ws.Cells.Item(1, 5).Value = "ResponseIndex"
$ws.Cells.Item(1, 6).Value = "ReactionTime"
$ws.Cells.Item(1, 7).Value = "Condition"

# Row 2 - update reaction time value, keep rest
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = "Gosto de todo o tipo de jogos e passatempos."
$ws.Cells.Item(2, 4).Value = "Completamente Verdadeiro"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.89754470000116271
$ws.Cells.Item(2, 7).Value = "active"

# Row 3 - cleared data (reset to zero/empty)
$ws.Cells.Item(3, 1).Value = 0
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = ""
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = ""

# Row 4 - cleared data (reset to zero/empty)
$ws.Cells.Item(4, 1).Value = 0
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = ""
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = ""

# Row 5 - cleared D column (reset to zero/empty), C5 and G5 already empty
$ws.Cells.Item(5, 1).Value = 0
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = ""
$ws.Cells.Item(5, 4).Value = ""
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = ""

# Autofit columns to reflect new (shorter) content widths
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(6).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(7).EntireColumn.AutoFit() | Out-Null
